$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4 (groundcover / Tyler Drago): status moves from 2nd Interview to 3rd Interview
# and the action date moves from 45989 to 45992.
$ws.Range("E4").Value = "3rd Interview"
$ws.Range("F4").Value = 45992

# Row 8 (LanceDB / Gregory West, 2nd Interview, 45980) was removed entirely.
# Deleting the entire row shifts all subsequent rows up by one, which matches
# the rest of the diff (rows 9-14 each becoming the prior row's values).
$ws.Rows.Item(8).Delete()
